$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 105.3125
$ws.Cells.Item(9, 9).Value = 109.5
$ws.Cells.Item(9, 10).Value = 92.75
$ws.Cells.Item(9, 11).Value = 109.5
$ws.Cells.Item(9, 12).Value = 92.75
$ws.Cells.Item(9, 13).Value = 59.5
$ws.Cells.Item(9, 14).Value = -430.75
$ws.Cells.Item(19, 8).Value = 1044.3636
$ws.Cells.Item(19, 9).Value = 1132.6666
$ws.Cells.Item(19, 10).Value = 1011.25
$ws.Cells.Item(19, 11).Value = 1132.6666
$ws.Cells.Item(19, 12).Value = 1011.25
$ws.Cells.Item(19, 13).Value = -957.6666
$ws.Cells.Item(19, 14).Value = -1361.25
$ws.Cells.Item(70, 8).Value = 2992.5334
$ws.Cells.Item(70, 9).Value = 2890.5
$ws.Cells.Item(70, 10).Value = 3060.5557
$ws.Cells.Item(70, 11).Value = 8671.5
$ws.Cells.Item(70, 12).Value = 9181.667099999999
$ws.Cells.Item(70, 13).Value = -8401.5
$ws.Cells.Item(70, 14).Value = -9721.667099999999
$ws.Cells.Item(73, 8).Value = 2992.5334
$ws.Cells.Item(73, 9).Value = 2890.5
$ws.Cells.Item(73, 10).Value = 3060.5557
$ws.Cells.Item(73, 11).Value = 8671.5
$ws.Cells.Item(73, 12).Value = 9181.667099999999
$ws.Cells.Item(73, 13).Value = -7735.5
$ws.Cells.Item(73, 14).Value = -11053.6671
$ws.Cells.Item(80, 8).Value = 1756.7858
$ws.Cells.Item(80, 9).Value = 1400.7142
$ws.Cells.Item(80, 10).Value = 2112.8572
$ws.Cells.Item(80, 11).Value = 4202.142599999999
$ws.Cells.Item(80, 12).Value = 6338.571599999999
$ws.Cells.Item(80, 13).Value = -3204.142599999999
$ws.Cells.Item(80, 14).Value = -8334.571599999999
$ws.Cells.Item(83, 8).Value = 1756.7858
$ws.Cells.Item(83, 9).Value = 1400.7142
$ws.Cells.Item(83, 10).Value = 2112.8572
$ws.Cells.Item(83, 11).Value = 12606.4278
$ws.Cells.Item(83, 12).Value = 19015.7148
$ws.Cells.Item(83, 13).Value = -7614.427799999999
$ws.Cells.Item(83, 14).Value = -28999.7148
$ws.Cells.Item(86, 8).Value = 17883.7
$ws.Cells.Item(86, 9).Value = 16104.625
$ws.Cells.Item(86, 10).Value = 25000
$ws.Cells.Item(86, 11).Value = 16104.625
$ws.Cells.Item(86, 12).Value = 25000
$ws.Cells.Item(86, 13).Value = -14981.625
$ws.Cells.Item(86, 14).Value = -27246
$ws.Cells.Item(89, 8).Value = 17883.7
$ws.Cells.Item(89, 9).Value = 16104.625
$ws.Cells.Item(89, 10).Value = 25000
$ws.Cells.Item(89, 11).Value = 80523.125
$ws.Cells.Item(89, 12).Value = 125000
$ws.Cells.Item(89, 13).Value = -74907.125
$ws.Cells.Item(89, 14).Value = -136232
$ws.Cells.Item(93, 8).Value = 34999
$ws.Cells.Item(93, 9).Value = 34999
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 34999
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = -32503
$ws.Cells.Item(93, 14).Value = ""
$ws.Cells.Item(98, 8).Value = 1011.36365
$ws.Cells.Item(98, 9).Value = 855.2778
$ws.Cells.Item(98, 10).Value = 1713.75
$ws.Cells.Item(98, 11).Value = 855.2778
$ws.Cells.Item(98, 12).Value = 1713.75
$ws.Cells.Item(98, 13).Value = 642.7222
$ws.Cells.Item(98, 14).Value = -4709.75
$ws.Cells.Item(122, 8).Value = 1011.36365
$ws.Cells.Item(122, 9).Value = 855.2778
$ws.Cells.Item(122, 10).Value = 1713.75
$ws.Cells.Item(122, 11).Value = 2565.8334
$ws.Cells.Item(122, 12).Value = 5141.25
$ws.Cells.Item(122, 13).Value = -115.8334
$ws.Cells.Item(122, 14).Value = -10041.25
$ws.Cells.Item(131, 8).Value = 1586.75
$ws.Cells.Item(131, 9).Value = 1799
$ws.Cells.Item(131, 10).Value = 950
$ws.Cells.Item(131, 11).Value = 5397
$ws.Cells.Item(131, 12).Value = 2850
$ws.Cells.Item(131, 13).Value = -357
$ws.Cells.Item(131, 14).Value = -12930
$ws.Cells.Item(132, 8).Value = 9387.956
$ws.Cells.Item(132, 9).Value = 11724.944
$ws.Cells.Item(132, 10).Value = 974.8
$ws.Cells.Item(132, 11).Value = 35174.83199999999
$ws.Cells.Item(132, 12).Value = 2924.4
$ws.Cells.Item(132, 13).Value = -32644.83199999999
$ws.Cells.Item(132, 14).Value = -7984.4
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(44, 8).Value = 29999
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 29999
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 12).Value = 29999
$ws.Cells.Item(44, 14).Value = -30975
$ws.Cells.Item(45, 8).Value = 3050.6667
$ws.Cells.Item(45, 9).Value = 3050.6667
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 11).Value = 3050.6667
$ws.Cells.Item(45, 12).Value = 0
$ws.Cells.Item(45, 13).Value = -2673.6667
$ws.Cells.Item(74, 8).Value = 7002.7856
$ws.Cells.Item(74, 9).Value = 3756.8823
$ws.Cells.Item(74, 10).Value = 20797.875
$ws.Cells.Item(74, 11).Value = 3756.8823
$ws.Cells.Item(74, 12).Value = 20797.875
$ws.Cells.Item(74, 13).Value = -2882.8823
$ws.Cells.Item(74, 14).Value = -22545.875
$ws.Cells.Item(77, 8).Value = 7002.7856
$ws.Cells.Item(77, 9).Value = 3756.8823
$ws.Cells.Item(77, 10).Value = 20797.875
$ws.Cells.Item(77, 11).Value = 18784.4115
$ws.Cells.Item(77, 12).Value = 103989.375
$ws.Cells.Item(77, 13).Value = -14416.4115
$ws.Cells.Item(77, 14).Value = -112725.375
$ws.Cells.Item(102, 8).Value = 1633.5714
$ws.Cells.Item(102, 9).Value = 1572.5
$ws.Cells.Item(102, 10).Value = 2000
$ws.Cells.Item(102, 11).Value = 1572.5
$ws.Cells.Item(102, 12).Value = 2000
$ws.Cells.Item(102, 13).Value = 49.5
$ws.Cells.Item(102, 14).Value = -5244
$ws.Cells.Item(132, 8).Value = 4795.4653
$ws.Cells.Item(132, 9).Value = 2735.1462
$ws.Cells.Item(132, 10).Value = 9764.471
$ws.Cells.Item(132, 11).Value = 8205.438600000001
$ws.Cells.Item(132, 12).Value = 29293.413
$ws.Cells.Item(132, 13).Value = -5675.438600000001
$ws.Cells.Item(132, 14).Value = -34353.413
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 0
$ws.Cells.Item(105, 9).Value = 0
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 0
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = ""
$ws.Cells.Item(105, 14).Value = ""
$ws.Cells.Item(107, 8).Value = 17858238
$ws.Cells.Item(107, 9).Value = 20834404
$ws.Cells.Item(107, 10).Value = 1240
$ws.Cells.Item(107, 11).Value = 20834404
$ws.Cells.Item(107, 12).Value = 1240
$ws.Cells.Item(107, 13).Value = -20832484
$ws.Cells.Item(107, 14).Value = -5080
$ws.Cells.Item(134, 8).Value = 2635.8823
$ws.Cells.Item(134, 9).Value = 2363.125
$ws.Cells.Item(134, 10).Value = 7000
$ws.Cells.Item(134, 11).Value = 7089.375
$ws.Cells.Item(134, 12).Value = 21000
$ws.Cells.Item(134, 13).Value = -4554.375
$ws.Cells.Item(134, 14).Value = -26070
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 32024.5
$ws.Cells.Item(16, 9).Value = 40999.332
$ws.Cells.Item(16, 10).Value = 5100
$ws.Cells.Item(16, 11).Value = 40999.332
$ws.Cells.Item(16, 12).Value = 5100
$ws.Cells.Item(16, 13).Value = -40712.332
$ws.Cells.Item(16, 14).Value = -5674
$ws.Cells.Item(58, 8).Value = 3531.5625
$ws.Cells.Item(58, 9).Value = 2350.3704
$ws.Cells.Item(58, 10).Value = 9910
$ws.Cells.Item(58, 11).Value = 2350.3704
$ws.Cells.Item(58, 12).Value = 9910
$ws.Cells.Item(58, 13).Value = -2147.3704
$ws.Cells.Item(58, 14).Value = -10316
$ws.Cells.Item(60, 8).Value = 26499.666
$ws.Cells.Item(60, 9).Value = 0
$ws.Cells.Item(60, 10).Value = 26499.666
$ws.Cells.Item(60, 11).Value = 0
$ws.Cells.Item(60, 12).Value = 26499.666
$ws.Cells.Item(60, 14).Value = -27521.666
$ws.Cells.Item(94, 8).Value = 6290.6313
$ws.Cells.Item(94, 9).Value = 12918.25
$ws.Cells.Item(94, 10).Value = 1470.5454
$ws.Cells.Item(94, 11).Value = 12918.25
$ws.Cells.Item(94, 12).Value = 1470.5454
$ws.Cells.Item(94, 13).Value = -12467.25
$ws.Cells.Item(94, 14).Value = -2372.5454
$ws.Cells.Item(105, 8).Value = 7881.4116
$ws.Cells.Item(105, 9).Value = 8265.6
$ws.Cells.Item(105, 10).Value = 5000
$ws.Cells.Item(105, 11).Value = 8265.6
$ws.Cells.Item(105, 12).Value = 5000
$ws.Cells.Item(105, 13).Value = -6518.6
$ws.Cells.Item(105, 14).Value = -8494
$ws.Cells.Item(113, 8).Value = 32024.5
$ws.Cells.Item(113, 9).Value = 40999.332
$ws.Cells.Item(113, 10).Value = 5100
$ws.Cells.Item(113, 11).Value = 40999.332
$ws.Cells.Item(113, 12).Value = 5100
$ws.Cells.Item(113, 13).Value = -38829.332
$ws.Cells.Item(113, 14).Value = -9440
$ws.Cells.Item(132, 8).Value = 2275.75
$ws.Cells.Item(132, 9).Value = 2305.4187
$ws.Cells.Item(132, 10).Value = 1000
$ws.Cells.Item(132, 11).Value = 6916.256100000001
$ws.Cells.Item(132, 12).Value = 3000
$ws.Cells.Item(132, 13).Value = -4386.256100000001
$ws.Cells.Item(132, 14).Value = -8060
$ws.Cells.Item(136, 8).Value = 3531.5625
$ws.Cells.Item(136, 9).Value = 2350.3704
$ws.Cells.Item(136, 10).Value = 9910
$ws.Cells.Item(136, 11).Value = 7051.111199999999
$ws.Cells.Item(136, 12).Value = 29730
$ws.Cells.Item(136, 13).Value = -4501.111199999999
$ws.Cells.Item(136, 14).Value = -34830
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 98.55556
$ws.Cells.Item(14, 9).Value = 98.55556
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 295.66668
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = -122.66668
$ws.Cells.Item(34, 8).Value = 2075.75
$ws.Cells.Item(34, 9).Value = 2565
$ws.Cells.Item(34, 10).Value = 1977.9
$ws.Cells.Item(34, 11).Value = 7695
$ws.Cells.Item(34, 12).Value = 5933.700000000001
$ws.Cells.Item(34, 13).Value = -7611
$ws.Cells.Item(34, 14).Value = -6101.700000000001
$ws.Cells.Item(59, 9).Value = 0
$ws.Cells.Item(59, 10).Value = 1500
$ws.Cells.Item(59, 11).Value = 0
$ws.Cells.Item(59, 12).Value = 4500
$ws.Cells.Item(59, 14).Value = -5580
$ws.Cells.Item(59, 13).Value = ""
$ws.Cells.Item(122, 8).Value = 1614853.6
$ws.Cells.Item(122, 9).Value = 16129032
$ws.Cells.Item(122, 10).Value = 2167.111
$ws.Cells.Item(122, 11).Value = 145161288
$ws.Cells.Item(122, 12).Value = 19503.999
$ws.Cells.Item(122, 13).Value = -145158838
$ws.Cells.Item(122, 14).Value = -24403.999
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2175.5
$ws.Cells.Item(102, 9).Value = 2719.6667
$ws.Cells.Item(102, 10).Value = 1196
$ws.Cells.Item(102, 11).Value = 2719.6667
$ws.Cells.Item(102, 12).Value = 1196
$ws.Cells.Item(102, 13).Value = -1097.6667
$ws.Cells.Item(102, 14).Value = -4440
$ws.Cells.Item(113, 8).Value = 1734.6086
$ws.Cells.Item(113, 9).Value = 1718.8572
$ws.Cells.Item(113, 10).Value = 1900
$ws.Cells.Item(113, 11).Value = 1718.8572
$ws.Cells.Item(113, 12).Value = 1900
$ws.Cells.Item(113, 13).Value = 451.1428000000001
$ws.Cells.Item(113, 14).Value = -6240
$ws.Cells.Item(122, 8).Value = 3322.1667
$ws.Cells.Item(122, 9).Value = 3178.353
$ws.Cells.Item(122, 10).Value = 3671.4285
$ws.Cells.Item(122, 11).Value = 9535.059000000001
$ws.Cells.Item(122, 12).Value = 11014.2855
$ws.Cells.Item(122, 13).Value = -7085.059000000001
$ws.Cells.Item(122, 14).Value = -15914.2855
$ws.Cells.Item(126, 8).Value = 2691.4614
$ws.Cells.Item(126, 9).Value = 2327.2856
$ws.Cells.Item(126, 10).Value = 3116.3333
$ws.Cells.Item(126, 11).Value = 6981.8568
$ws.Cells.Item(126, 12).Value = 9348.999899999999
$ws.Cells.Item(126, 13).Value = -4511.8568
$ws.Cells.Item(126, 14).Value = -14288.9999
$ws.Cells.Item(132, 8).Value = 17036.857
$ws.Cells.Item(132, 9).Value = 17036.857
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 51110.571
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -48580.571
$ws.Cells.Item(132, 14).Value = ""
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5266.4443
$ws.Cells.Item(7, 9).Value = 5199.857
$ws.Cells.Item(7, 10).Value = 5499.5
$ws.Cells.Item(7, 11).Value = 5199.857
$ws.Cells.Item(7, 12).Value = 5499.5
$ws.Cells.Item(7, 13).Value = -5087.857
$ws.Cells.Item(7, 14).Value = -5723.5
$ws.Cells.Item(68, 8).Value = 16000.667
$ws.Cells.Item(68, 9).Value = 10000
$ws.Cells.Item(68, 10).Value = 19001
$ws.Cells.Item(68, 11).Value = 10000
$ws.Cells.Item(68, 12).Value = 19001
$ws.Cells.Item(68, 13).Value = -9251
$ws.Cells.Item(68, 14).Value = -20499
$ws.Cells.Item(71, 8).Value = 16000.667
$ws.Cells.Item(71, 9).Value = 10000
$ws.Cells.Item(71, 10).Value = 19001
$ws.Cells.Item(71, 11).Value = 50000
$ws.Cells.Item(71, 12).Value = 95005
$ws.Cells.Item(71, 13).Value = -46256
$ws.Cells.Item(71, 14).Value = -102493
$ws.Cells.Item(126, 8).Value = 5266.4443
$ws.Cells.Item(126, 9).Value = 5199.857
$ws.Cells.Item(126, 10).Value = 5499.5
$ws.Cells.Item(126, 11).Value = 15599.571
$ws.Cells.Item(126, 12).Value = 16498.5
$ws.Cells.Item(126, 13).Value = -13129.571
$ws.Cells.Item(126, 14).Value = -21438.5
$ws.Cells.Item(136, 8).Value = 1535.7894
$ws.Cells.Item(136, 9).Value = 1343.3334
$ws.Cells.Item(136, 10).Value = 5000
$ws.Cells.Item(136, 11).Value = 4030.0002
$ws.Cells.Item(136, 12).Value = 15000
$ws.Cells.Item(136, 13).Value = -1480.0002
$ws.Cells.Item(136, 14).Value = -20100
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 20013.691
$ws.Cells.Item(41, 9).Value = 10000
$ws.Cells.Item(41, 10).Value = 21834.363
$ws.Cells.Item(41, 11).Value = 10000
$ws.Cells.Item(41, 12).Value = 21834.363
$ws.Cells.Item(41, 13).Value = -9610
$ws.Cells.Item(41, 14).Value = -22614.363
$ws.Cells.Item(122, 8).Value = 42868.215
$ws.Cells.Item(122, 9).Value = 2737.8572
$ws.Cells.Item(122, 10).Value = 163259.28
$ws.Cells.Item(122, 11).Value = 8213.571599999999
$ws.Cells.Item(122, 12).Value = 489777.84
$ws.Cells.Item(122, 13).Value = -5763.571599999999
$ws.Cells.Item(122, 14).Value = -494677.84
$ws.Cells.Item(126, 8).Value = 1905.7368
$ws.Cells.Item(126, 9).Value = 1605.6364
$ws.Cells.Item(126, 10).Value = 2318.375
$ws.Cells.Item(126, 11).Value = 4816.9092
$ws.Cells.Item(126, 12).Value = 6955.125
$ws.Cells.Item(126, 13).Value = -2346.9092
$ws.Cells.Item(126, 14).Value = -11895.125
$ws.Cells.Item(132, 8).Value = 2742.1875
$ws.Cells.Item(132, 9).Value = 2118.7812
$ws.Cells.Item(132, 10).Value = 3989
$ws.Cells.Item(132, 11).Value = 6356.3436
$ws.Cells.Item(132, 12).Value = 11967
$ws.Cells.Item(132, 13).Value = -3826.3436
$ws.Cells.Item(132, 14).Value = -17027
